$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Capture current (pre-edit) values for row 2 and row 3 that need to be swapped ---
# (use Value2 for reading - Value getter is unreliable in this environment)
$A2 = $ws.Range("A2").Value2
$B2 = $ws.Range("B2").Value2
$E2 = $ws.Range("E2").Value2
$F2 = $ws.Range("F2").Value2
$G2 = $ws.Range("G2").Value2
$H2 = $ws.Range("H2").Value2

$A3 = $ws.Range("A3").Value2
$B3 = $ws.Range("B3").Value2
$E3 = $ws.Range("E3").Value2
$F3 = $ws.Range("F3").Value2
$G3 = $ws.Range("G3").Value2
$H3 = $ws.Range("H3").Value2

# --- Swap row 2 <-> row 3 for columns A, B, E, F, G, H ---
$ws.Range("A2").Value = $A3
$ws.Range("B2").Value = $B3
$ws.Range("E2").Value = $E3
$ws.Range("F2").Value = $F3
$ws.Range("G2").Value = $G3
$ws.Range("H2").Value = $H3

$ws.Range("A3").Value = $A2
$ws.Range("B3").Value = $B2
$ws.Range("E3").Value = $E2
$ws.Range("F3").Value = $F2
$ws.Range("G3").Value = $G2
$ws.Range("H3").Value = $H2

# --- Update Ost (Q) / Nord (R) coordinates: the new values are the rounded versions
#     of the coordinates belonging to the record that now occupies that row ---
$ws.Range("Q2").Value = 431104
$ws.Range("R2").Value = 6811804

$ws.Range("Q3").Value = 431106
$ws.Range("R3").Value = 6811802

$ws.Range("Q4").Value = 431104
$ws.Range("R4").Value = 6811805

# --- Remove the Starttid (Z) and Sluttid (AB) cells for rows 2, 3 and 4 ---
$ws.Range("Z2").ClearContents()
$ws.Range("AB2").ClearContents()

$ws.Range("Z3").ClearContents()
$ws.Range("AB3").ClearContents()

$ws.Range("Z4").ClearContents()
$ws.Range("AB4").ClearContents()
